# Auto-generated Excel COM-interop edit script
# Applies the "Updated symbol list" commit changes to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell "D2" "235.73"
Set-TextCell "G2" "20"
Set-TextCell "D3" "21.71"
Set-TextCell "G3" "20"
Set-TextCell "B4" "HuobiToken"
Set-TextCell "C4" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D4" "5.361"
Set-TextCell "E4" "3HuobiTokenHT"
Set-TextCell "G4" "20"
Set-TextCell "B5" "Cronos"
Set-TextCell "C5" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D5" "0.05569"
Set-TextCell "E5" "4CronosCRO"
Set-TextCell "G5" "20"
Set-TextCell "B6" "GateToken"
Set-TextCell "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D6" "3.366"
Set-TextCell "E6" "5GateTokenGT"
Set-TextCell "G6" "20"
Set-TextCell "B7" "KuCoinToken"
Set-TextCell "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell "D7" "6.458"
Set-TextCell "E7" "6KuCoinTokenKCS"
Set-TextCell "G7" "20"
Set-TextCell "B8" "MXToken"
Set-TextCell "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D8" "0.8007"
Set-TextCell "E8" "7MXTokenMX"
Set-TextCell "G8" "20"
Set-TextCell "B9" "FTXToken"
Set-TextCell "C9" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D9" "1.040"
Set-TextCell "E9" "8FTXTokenFTT"
Set-TextCell "G9" "20"
Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1403"
Set-TextCell "E10" "9WazirXWRX"
Set-TextCell "G10" "20"
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.07241"
Set-TextCell "E11" "10MandalaExchangeTokenMDX"
Set-TextCell "G11" "20"
Set-TextCell "B12" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D12" "0.03179"
Set-TextCell "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextCell "G12" "20"
Set-TextCell "B13" "BitrueCoin"
Set-TextCell "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.02933"
Set-TextCell "E13" "12BitrueCoinBTR"
Set-TextCell "G13" "20"
Set-TextCell "B14" "BitMartToken"
Set-TextCell "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.09237"
Set-TextCell "E14" "13BitMartTokenBMX"
Set-TextCell "G14" "20"
Set-TextCell "B15" "BitForexToken"
Set-TextCell "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001658"
Set-TextCell "E15" "14BitForexTokenBF"
Set-TextCell "G15" "20"
Set-TextCell "B16" "MCDex"
Set-TextCell "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D16" "3.253"
Set-TextCell "E16" "15MCDexMCB"
Set-TextCell "G16" "20"
Set-TextCell "B17" "CoinExToken"
Set-TextCell "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D17" "0.04755"
Set-TextCell "E17" "16CoinExTokenCET"
Set-TextCell "G17" "20"
Set-TextCell "B18" "One"
Set-TextCell "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D18" "0.0005710"
Set-TextCell "E18" "17OneONE"
Set-TextCell "G18" "20"
Set-TextCell "B19" "TigerCash"
Set-TextCell "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D19" "0.006257"
Set-TextCell "E19" "18TigerCashTCH"
Set-TextCell "G19" "20"
Set-TextCell "B20" "HotbitToken"
Set-TextCell "C20" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell "D20" "0.005077"
Set-TextCell "E20" "19HotbitTokenHTB"
Set-TextCell "G20" "20"
Set-TextCell "B21" "BitKan"
Set-TextCell "C21" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell "D21" "0.001050"
Set-TextCell "E21" "20BitKanKAN"
Set-TextCell "G21" "20"
Set-TextCell "B22" "NitroEx"
Set-TextCell "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell "D22" "0.0001500"
Set-TextCell "E22" "21NitroExNTX"
Set-TextCell "G22" "20"
Set-TextCell "B23" "UpBots"
Set-TextCell "C23" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextCell "D23" "0.0004200"
Set-TextCell "E23" "22UpBotsUBXT"
Set-TextCell "G23" "20"
Set-TextCell "B24" "LEO"
Set-TextCell "C24" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D24" "3.936"
Set-TextCell "E24" "23LEOLEO"
Set-TextCell "G24" "20"
Set-TextCell "G25" "20"
Set-TextCell "G26" "20"
Set-TextCell "G27" "20"
Set-TextCell "G28" "20"
Set-TextCell "G29" "20"
Set-TextCell "G30" "20"
Set-TextCell "G31" "20"
Set-TextCell "G32" "20"
Set-TextCell "G33" "20"
Set-TextCell "G34" "20"
Set-TextCell "G35" "20"
Set-TextCell "G36" "20"
Set-TextCell "G37" "20"
Set-TextCell "G38" "20"
Set-TextCell "G39" "20"
Set-TextCell "D40" "0.04117"
Set-TextCell "G40" "20"
Set-TextCell "D41" "0.007037"
Set-TextCell "G41" "20"
Set-TextCell "D42" "0.003500"
Set-TextCell "E42" "41CEJICEJIBestin24h"
Set-TextCell "G42" "20"
Set-TextCell "D43" "0.1038"
Set-TextCell "G43" "20"
Set-TextCell "D44" "0.008868"
Set-TextCell "G44" "20"
Set-TextCell "D45" "0.00005433"
Set-TextCell "G45" "20"
Set-TextCell "G46" "20"
Set-TextCell "D47" "0.6799"
Set-TextCell "G47" "20"
Set-TextCell "D48" "0.03299"
Set-TextCell "G48" "20"
Set-TextCell "D49" "0.00002100"
Set-TextCell "G49" "20"
Set-TextCell "D50" "0.01010"
Set-TextCell "G50" "20"
Set-TextCell "G51" "20"
